$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-05 Monday" "2025-05-06 Tuesday"

Replace-Text "876×7=6132" "866×3=2598"
Replace-Text "618×7=4326" "930×4=3720"
Replace-Text "209×8=1672" "947×6=5682"
Replace-Text "606×4=2424" "331×8=2648"
Replace-Text "106×7=742" "536×6=3216"

Replace-Text "733×8=5864" "684×9=6156"
Replace-Text "825×3=2475" "946×8=7568"
Replace-Text "999×9=8991" "123×3=369"
Replace-Text "628×5=3140" "812×8=6496"
Replace-Text "640×5=3200" "692×5=3460"

Replace-Text "237×9=2133" "166×8=1328"
Replace-Text "138×3=414" "743×7=5201"
Replace-Text "501×2=1002" "494×8=3952"
Replace-Text "213×9=1917" "423×2=846"
Replace-Text "775×9=6975" "366×3=1098"

Replace-Text "163×7=1141" "489×4=1956"
Replace-Text "819×6=4914" "196×6=1176"
Replace-Text "609×7=4263" "868×8=6944"
Replace-Text "506×3=1518" "578×9=5202"
Replace-Text "187×9=1683" "210×7=1470"

Replace-Text "864×4=3456" "868×8=6944"
Replace-Text "779×3=2337" "723×9=6507"
Replace-Text "514×4=2056" "702×4=2808"
Replace-Text "317×7=2219" "150×6=900"
Replace-Text "335×3=1005" "430×3=1290"
